# Update the public EPEX Spot prices workbook with the latest day of data.

# Helper: write a plain-text date-looking value (e.g. "2025-08-02") into a
# cell without letting Excel auto-convert it to a real date serial number,
# and without leaving a non-default number format applied to the cell
# (matches the existing "YYYY-MM-DD" text cells already in the sheet).
function Set-DateText($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column (AZ) with the 04-aug prices.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header style (bold + border + centered) from the previous day's
# header cell (AY1) onto the new header cell (AZ1).
$wsPrix.Cells.Item(1, 51).Copy() | Out-Null
$wsPrix.Cells.Item(1, 52).PasteSpecial(-4122) | Out-Null
$wsPrix.Cells.Item(1, 52).Value = "04-aug"

$prixValues = @(
    80.98999999999999,
    70.06999999999999,
    69.97,
    63.8,
    61.78,
    65.7,
    74.33,
    78.83,
    78.43000000000001,
    52.57,
    30,
    13.75,
    6.76,
    4.05,
    0.65,
    0,
    5.79,
    7.28,
    38.94,
    63.7,
    70.25,
    74,
    67.14,
    49.4
)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 52).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the two missing days (2025-08-02 and 2025-08-03).
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

Set-DateText $wsGaz.Cells.Item(49, 1) "2025-08-02"
$wsGaz.Cells.Item(49, 2).Value = 32.775

Set-DateText $wsGaz.Cells.Item(50, 1) "2025-08-03"
$wsGaz.Cells.Item(50, 2).Value = 32.775

# ---------------------------------------------------------------------------
# Sheet "CO2": append the two missing days (2025-08-02 and 2025-08-03).
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

Set-DateText $wsCo2.Cells.Item(49, 1) "2025-08-02"
$wsCo2.Cells.Item(49, 2).Value = 70.58

Set-DateText $wsCo2.Cells.Item(50, 1) "2025-08-03"
$wsCo2.Cells.Item(50, 2).Value = 70.58
